# Re-apply the latest scraped crypto market snapshot to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.312.53"
$ws.Range("E2").Value = "  -2.96%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.217.28"
$ws.Range("E3").Value = "  -2.32%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "108.27"
$ws.Range("E5").Value = "  -10.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "296.05"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").Value = "  -3.28%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.596"
$ws.Range("E9").Value = "  -4.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.58"
$ws.Range("E10").Value = "  -8.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0910"
$ws.Range("E11").Value = "  -3.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.50"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.80"
$ws.Range("E13").Value = "  -5.59%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.979"
$ws.Range("E14").Value = "  +7.48%  "
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.103"
$ws.Range("E15").Value = "  -2.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.98"
$ws.Range("E16").Value = "  -2.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.548.56"
$ws.Range("E17").Value = "  -2.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.224.68"
$ws.Range("E18").Value = "  -1.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.297.88"
$ws.Range("E19").Value = "  -2.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.39"
$ws.Range("E20").Value = "  +7.04%  "
$ws.Range("E21").Value = "  -4.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.37"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.48"
$ws.Range("E23").Value = "  +21.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.31"
$ws.Range("E24").Value = "  -3.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "228.12"
$ws.Range("E25").Value = "  -2.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.00"
$ws.Range("E26").Value = "  -5.55%  "
$ws.Range("E27").Value = "  -1.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.60"
$ws.Range("E28").Value = "  -2.77%  "
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.45"
$ws.Range("E30").Value = "  -8.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.21"
$ws.Range("E31").Value = "  -4.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "173.60"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.92"
$ws.Range("E33").Value = "  -3.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0894"
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.59"
$ws.Range("E35").Value = "  -2.50%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.04"
$ws.Range("E36").Value = "  +10.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.38"
$ws.Range("E37").Value = "  +1.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.126"
$ws.Range("E38").Value = "  -3.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0366"
$ws.Range("E39").Value = "  -2.68%  "
$ws.Range("E40").Value = "  -4.06%  "
$ws.Range("E41").Value = "  -4.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.28"
$ws.Range("E42").Value = "  -2.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.232"
$ws.Range("E43").Value = "  -2.43%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.58"
$ws.Range("E45").Value = "  -9.48%  "
$ws.Range("E46").Value = "  -5.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.41"
$ws.Range("E47").Value = "  -6.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.32"
$ws.Range("E48").Value = "  +4.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.50"
$ws.Range("E49").Value = "  +2.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.41"
$ws.Range("E50").Value = "  -1.23%  "
$ws.Range("E51").Value = "  +4.29%  "
